$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.529.49"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "3.785.11"
$ws.Range("E3").Value = "  +0.51%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'614.32"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("D6").Value = "'177.31"
$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("D7").Value = "3.781.05"
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("D11").Value = "'6.42"
$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("D12").Value = "'0.483"
$ws.Range("E12").Value = "  -1.85%  "

$ws.Range("D13").Value = "'39.82"
$ws.Range("E13").Value = "  -3.81%  "

$ws.Range("D14").Value = "'0.0000254"
$ws.Range("E14").Value = "  -2.68%  "

$ws.Range("D15").Value = "4.418.56"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").Value = "3.787.67"

$ws.Range("D17").Value = "69.625.45"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "'7.55"
$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("E19").Value = "  -3.69%  "

$ws.Range("D20").Value = "'509.07"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").Value = "'16.56"
$ws.Range("E21").Value = "  -1.39%  "

$ws.Range("D22").Value = "'9.59"
$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("D23").Value = "'0.734"
$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("D24").Value = "'2.47"
$ws.Range("E24").Value = "  -1.92%  "

$ws.Range("D25").Value = "'86.27"
$ws.Range("E25").Value = "  -1.33%  "

$ws.Range("D26").Value = "'12.86"
$ws.Range("E26").Value = "  -2.69%  "

$ws.Range("D27").Value = "'0.0000141"
$ws.Range("E27").Value = "  +2.70%  "

$ws.Range("D28").Value = "'10.55"

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "'3.00"
$ws.Range("E30").Value = "  +2.95%  "

$ws.Range("E31").Value = "  -0.54%  "

$ws.Range("D32").Value = "'8.11"
$ws.Range("E32").Value = "  +2.41%  "

$ws.Range("D33").Value = "'31.29"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("E36").Value = "  -1.33%  "

$ws.Range("D37").Value = "'6.13"
$ws.Range("E37").Value = "  -1.68%  "

$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "  +6.05%  "

$ws.Range("D39").Value = "'480.23"
$ws.Range("E39").Value = "  +11.51%  "

$ws.Range("D40").Value = "'0.340"
$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").Value = "'2.06"
$ws.Range("E41").Value = "  -3.04%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'49.77"
$ws.Range("E42").Value = "  -0.95%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.98"
$ws.Range("E43").Value = "  +4.16%  "

$ws.Range("D44").Value = "'44.04"
$ws.Range("E44").Value = "  -3.68%  "

$ws.Range("D45").Value = "'8.57"
$ws.Range("E45").Value = "  -2.20%  "

$ws.Range("D46").Value = "2.946.08"
$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("D47").Value = "'0.0364"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("D48").Value = "'27.44"
$ws.Range("E48").Value = "  -0.48%  "

$ws.Range("D49").Value = "'139.80"
$ws.Range("E49").Value = "  +2.20%  "

$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").Value = "'2.45"
$ws.Range("E51").Value = "  -2.30%  "
